$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.152.26'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.39'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.66'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4479'
$ws.Range("E7").Value = '  +5.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3761'
$ws.Range("E8").Value = '  +2.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07517'
$ws.Range("E9").Value = '  +3.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8766'
$ws.Range("E10").Value = '  +3.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.96'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.815.69'
$ws.Range("E12").Value = '  -0.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.755'
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.59'
$ws.Range("E14").Value = '  +5.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.363'
$ws.Range("E15").Value = '  +1.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07112'
$ws.Range("E16").Value = '  +1.08%  '

$ws.Range("E17").Value = '  -0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008778'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.06'
$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.175.87'
$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.246'
$ws.Range("E22").Value = '  +2.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.040.39'
$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.991'
$ws.Range("E25").Value = '  +0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.446'
$ws.Range("E26").Value = '  +8.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.89'
$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.54'
$ws.Range("E28").Value = '  +1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.360'
$ws.Range("E29").Value = '  +2.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.32'
$ws.Range("E30").Value = '  +1.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08824'
$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7708'
$ws.Range("E32").Value = '  +4.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.187'
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("E34").Value = '  +2.92%  '

$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9994'
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  +1.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01992'
$ws.Range("E38").Value = '  +2.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5348'
$ws.Range("E41").Value = '  +5.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1729'
$ws.Range("E42").Value = '  +2.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.855'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.248'
$ws.Range("E44").Value = '  +13.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.779'
$ws.Range("E45").Value = '  +2.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5119'
$ws.Range("E46").Value = '  +8.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.63'
$ws.Range("E47").Value = '  +1.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.708'
$ws.Range("E48").Value = '  +3.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.92'
$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9992'
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06365'
$ws.Range("E51").Value = '  +0.68%  '

# Row 39 and 40 content swap (FraxShare/Hedera order changed)
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05300'
$ws.Range("E39").Value = '  +1.22%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.439'
$ws.Range("E40").Value = '  +1.47%  '
